$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend table with a new column F (mirrors column E header, adds a new
# blank data cell in row 2 and a new JSON payload value in row 3).

# Copy the existing style (border + centered alignment, style index 1)
# from A1 onto the new F1:F3 cells before setting their values, so we
# don't introduce new/duplicate style records.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:F3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# F1 header mirrors E1 ("open")
$ws.Range("F1").Value = $ws.Range("E1").Text

# F2 stays blank (just inherits the style copied above)

# F3 holds the new JSON payload referencing ${url}
$ws.Range("F3").Value = '{"target":"${url}"}'

# Size the new column similarly to the existing bestFit columns
$ws.Columns("F").ColumnWidth = 16.14

# Move the active selection to F8, as in the target workbook
$ws.Range("F8").Select() | Out-Null
